$d = $word.ActiveDocument

# --- Step 1: merge the split runs / drop the proofErr spell-check markers ---
# Paragraphs 1, 2, 3, 6 and 9 had their text split across two runs with a
# <w:proofErr w:type="spellStart"/> ... <w:proofErr w:type="spellEnd"/> pair
# in between. Running a Find/Replace across the whole body (old text == new
# text) makes Word rebuild each paragraph as a single run and drop the
# proofErr markers that fall inside the replaced range.
$allText = "1 – Salto Abalacob`r" + `
           "2 – Salto cmj`r" + `
           "3 – Salto sj`r" + `
           "4 – Salto Continuo`r" + `
           "5 – Peso muerto`r" + `
           "6 – Velocidad 10 mts`r" + `
           "7 – Remo`r" + `
           "8 – Yoyo Test`r" + `
           "9 – Sentadilla Bulgara`r"

[void]$d.Content.Find.Execute($allText, $true, $false, $false, $false, $false, $true, 1, $false, $allText, 2)

# That replace leaves a stray (now run-less) proofErr marker trailing after
# paragraph 9, pushed into a brand-new, still-empty paragraph 10 — which is
# exactly the new line the commit adds, so we simply fill it in.
$p10 = $d.Paragraphs(10).Range
$p10.Text = "10 – Peso muerto 1 pierna"

# --- Step 2: clean up the stray proofErr marker trailing paragraph 10 ---
# Re-run the same trick across the now-complete text (including the new
# paragraph 10, with no trailing paragraph mark since it is the last
# paragraph) so the leftover proofErr marker has nowhere left to go and is
# dropped for good.
$allText2 = $allText + "10 – Peso muerto 1 pierna"
[void]$d.Content.Find.Execute($allText2, $true, $false, $false, $false, $false, $true, 1, $false, $allText2, 2)

# --- Step 3: restore the _GoBack bookmark on the new last paragraph ---
$p10 = $d.Paragraphs(10).Range
$bmPos = $p10.End - 1
$bmRange = $d.Range($bmPos, $bmPos)
[void]$d.Bookmarks.Add("_GoBack", $bmRange)
